# Updates for Moz training
# Replace accented login names in column A with their unaccented equivalents,
# and replace the "angelo" / "Ângelo Intimane" user (row 20) with the new
# "diquissone" / "Diquissone" user.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value  = "fetima"
$ws.Range("A18").Value = "denio"
$ws.Range("A20").Value = "diquissone"
$ws.Range("E20").Value = "Diquissone"
$ws.Range("A23").Value = "arsenio"
$ws.Range("A27").Value = "joao"
$ws.Range("A31").Value = "silvia"
$ws.Range("A35").Value = "felix"
$ws.Range("A39").Value = "gloria"
$ws.Range("A40").Value = "elisio"
$ws.Range("A42").Value = "eusebio"
$ws.Range("A43").Value = "angelica"
$ws.Range("A48").Value = "enia"
$ws.Range("A52").Value = "tania"
$ws.Range("A54").Value = "graca "
$ws.Range("A55").Value = "mauricio"
$ws.Range("A57").Value = "temotio"
$ws.Range("A69").Value = "nortencio"
$ws.Range("A70").Value = "hortencia"
$ws.Range("A77").Value = "amancio"

# Update the sheet view to reflect the last-selected cell (A20), matching the
# saved workbook view state (drops the old topLeftCell/selection on A65).
$ws.Range("A20").Select() | Out-Null
